$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (pushes the old rows 59..64 down to 60..65)
$ws.Rows(59).Insert()

# Populate the new row 59 with the new weekly price-report record
$ws.Range("A59").Value = 4
$ws.Range("B59").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value = "Los Lagos"
$ws.Range("D59").Value = 44984
$ws.Range("E59").Value = 10
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100101
$ws.Range("H59").Value = "Berries"
$ws.Range("I59").Value = 100101001
$ws.Range("J59").Value = "Arándano (blue)"
$ws.Range("K59").Value = "Sin especificar"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 200
$ws.Range("N59").Value = 2000
$ws.Range("O59").Value = 2200
$ws.Range("P59").Value = 2100
$ws.Range("Q59").Value = "$/bandeja 2 kilos"
$ws.Range("R59").Value = "Provincia de Curicó"
$ws.Range("S59").Value = 1050
$ws.Range("T59").Value = 2

# Keep the date cell formatted like the rest of the date column
$ws.Range("D59").NumberFormat = $ws.Range("D60").NumberFormat
